# Remove the discontinued course rows from the 2024B cupos summary.
# Deleting from bottom to top keeps the remaining row numbers stable
# while each EntireRow delete is applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(24, 23, 10, 7)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
